$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = @{
    'D2' = '37.735.70'
    'E2' = '  -1.34%  '
    'D3' = '2.030.63'
    'E3' = '  -1.86%  '
    'E4' = '  -0.17%  '
    'D5' = '227.49'
    'E5' = '  -1.35%  '
    'D6' = '0.603'
    'E6' = '  -1.74%  '
    'E7' = '  -2.49%  '
    'E8' = '  +0.03%  '
    'D9' = '0.376'
    'E9' = '  -3.06%  '
    'D10' = '0.0826'
    'E10' = '  +2.13%  '
    'E11' = '  -0.29%  '
    'D12' = '14.59'
    'E12' = '  -2.05%  '
    'D14' = '21.00'
    'E14' = '  -1.68%  '
    'D15' = '0.765'
    'E15' = '  -0.18%  '
    'E16' = '  -2.54%  '
    'D17' = '2.041.00'
    'E17' = '  -1.88%  '
    'D18' = '37.689.21'
    'E18' = '  -1.25%  '
    'D19' = '69.48'
    'E19' = '  -1.11%  '
    'E20' = '  -6.39%  '
    'D21' = '0.0₃0825'
    'E21' = '  -1.75%  '
    'D22' = '223.59'
    'E22' = '  -1.10%  '
    'E23' = '  +0.01%  '
    'D24' = '2.36'
    'E24' = '  -3.05%  '
    'D25' = '2.26'
    'E25' = '  +0.25%  '
    'D26' = '168.11'
    'E26' = '  +0.60%  '
    'E27' = '  +1.02%  '
    'E28' = '  -3.63%  '
    'D29' = '18.77'
    'E29' = '  -1.18%  '
    'E30' = '  -4.64%  '
    'E31' = '  +0.90%  '
    'E32' = '  +7.97%  '
    'D33' = '4.38'
    'E33' = '  -4.11%  '
    'D34' = '0.0605'
    'E34' = '  -0.44%  '
    'D35' = '4.48'
    'E35' = '  -3.21%  '
    'D36' = '6.46'
    'E36' = '  +2.40%  '
    'D37' = '2.30'
    'E37' = '  -2.10%  '
    'D38' = '3.42'
    'E38' = '  +2.51%  '
    'E39' = '  -0.01%  '
    'D40' = '18.19'
    'E40' = '  +5.92%  '
    'D41' = '1.538.93'
    'E41' = '  +0.92%  '
    'E42' = '  -1.33%  '
    'D43' = '95.60'
    'E43' = '  -2.78%  '
    'E44' = '  -3.24%  '
    'D45' = '0.0907'
    'E45' = '  -2.31%  '
    'D46' = '4.07'
    'E46' = '  +0.86%  '
    'E47' = '  -2.89%  '
    'E48' = '  -1.98%  '
    'E49' = '  +0.27%  '
    'D50' = '7.11'
    'E50' = '  -0.82%  '
    'D51' = '2.219.55'
    'E51' = '  -1.85%  '
}

foreach ($addr in $cellUpdates.Keys) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $cellUpdates[$addr]
    $r.Style = "Normal"
}
